$wb = $excel.ActiveWorkbook

$oldGuid = "b8b8e772-cbc2-4df0-aa6b-272efc54f374"
$newGuid = "4024b33c-99ac-4e1c-bc7c-e478e2256397"

$newHoDate = "2016-09-03 07:01:35"
$newZhDate = "2016-09-03 07:01:30"

$newMd   = "$newGuid.md"
$newPath = "e2e\$newGuid.md"
$newZhXlf = "$newGuid.30e87f6ec4442c953bcfddf9e1087afbc8ccdbdd.zh-cn.xlf"
$newDeXlf = "$newGuid.30e87f6ec4442c953bcfddf9e1087afbc8ccdbdd.de-de.xlf"

# The external hyperlink target (unchanged by this edit) used by all three sheets.
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66b89f60b400b6e171908e5beb94c62add0e478b/e2e/$oldGuid.md"

# ---------------- Overview sheet ----------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("B2").Value = $newPath
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkUrl, "", "", $newPath)
$wsOverview.Range("G2").Value = $newHoDate

# ---------------- zh-cn sheet ----------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $hyperlinkUrl, "", "", $newMd)
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhDate

# ---------------- de-de sheet ----------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $hyperlinkUrl, "", "", $newMd)
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newHoDate
